# TC06_Canine_Filter_Diagnosis-PulmNeoplasm.xlsx regression update
# (icdc regression 1 to 16 and 29 to 35)
#
# The workbook's "startup" sheet lists, one row per tab (CasesTab /
# SamplesTab / FilesTab / StudyFilesTab), the Neo4j "query" text used
# by that tab (column B) and the web "StatQuery" text (column C).
#
# This change fixes two of the Cypher queries stored in column B:
#   * CasesTab query:  the `Weight (kg)` projection now rounds the
#     value the same way `Age` already does.
#   * FilesTab query:  the `samp` variable used in the optional MATCH
#     for the study/diagnosis chain is now explicitly typed as
#     `:sample`, matching the rest of the query.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowByLabel([string]$label) {
    $used = $ws.UsedRange
    for ($r = 1; $r -le $used.Rows.Count; $r++) {
        if ($ws.Cells.Item($r, 1).Value2 -eq $label) {
            return $r
        }
    }
    return -1
}

# --- CasesTab query: fix the Weight (kg) coalesce expression ---------------
$oldWeight = 'coalesce(demo.weight, '''') AS `Weight (kg)`,'
$newWeight = 'coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '''') AS `Weight (kg)`,'

$casesRow = Get-RowByLabel "CasesTab"
if ($casesRow -gt 0) {
    $cell = $ws.Cells.Item($casesRow, 2)
    $text = $cell.Value2
    if ($text.Contains($oldWeight)) {
        $cell.Value2 = $text.Replace($oldWeight, $newWeight)
    }
}

# --- FilesTab query: type the samp variable as :sample ----------------------
$oldSamp = '<-[*]-(samp)'
$newSamp = '<-[*]-(samp:sample)'

$filesRow = Get-RowByLabel "FilesTab"
if ($filesRow -gt 0) {
    $cell = $ws.Cells.Item($filesRow, 2)
    $text = $cell.Value2
    if ($text.Contains($oldSamp)) {
        $cell.Value2 = $text.Replace($oldSamp, $newSamp)
    }
}
